$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine the last used row based on column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Update column C ("Förändrad") for every data row (2..lastRow) from 45172 to 45175
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45175
